$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C and column D values (trip planner distance/cost data) to match
# the resized / recalculated figures from the commit "added trip planner to gui
# and also resized gui".
$ws.Cells.Item(1, 3).Value = 851
$ws.Cells.Item(1, 4).Value = 620
$ws.Cells.Item(2, 3).Value = 1124
$ws.Cells.Item(2, 4).Value = 373
$ws.Cells.Item(3, 3).Value = 1100
$ws.Cells.Item(3, 4).Value = 429
$ws.Cells.Item(4, 3).Value = 918
$ws.Cells.Item(4, 4).Value = 591
$ws.Cells.Item(5, 3).Value = 778
$ws.Cells.Item(5, 4).Value = 399
$ws.Cells.Item(6, 3).Value = 884
$ws.Cells.Item(6, 4).Value = 428
$ws.Cells.Item(7, 3).Value = 585
$ws.Cells.Item(7, 4).Value = 655
$ws.Cells.Item(8, 3).Value = 431
$ws.Cells.Item(8, 4).Value = 458
$ws.Cells.Item(9, 3).Value = 865
$ws.Cells.Item(9, 4).Value = 360
$ws.Cells.Item(10, 3).Value = 52
$ws.Cells.Item(10, 4).Value = 519
$ws.Cells.Item(11, 3).Value = 624
$ws.Cells.Item(11, 4).Value = 729
$ws.Cells.Item(12, 3).Value = 809
$ws.Cells.Item(12, 4).Value = 460
$ws.Cells.Item(13, 3).Value = 178
$ws.Cells.Item(13, 4).Value = 652
$ws.Cells.Item(14, 3).Value = 126
$ws.Cells.Item(14, 4).Value = 598
$ws.Cells.Item(15, 3).Value = 730
$ws.Cells.Item(15, 4).Value = 587
$ws.Cells.Item(16, 3).Value = 933
$ws.Cells.Item(16, 4).Value = 804
$ws.Cells.Item(17, 3).Value = 757
$ws.Cells.Item(17, 4).Value = 320
$ws.Cells.Item(18, 3).Value = 658
$ws.Cells.Item(18, 4).Value = 297
$ws.Cells.Item(19, 3).Value = 720
$ws.Cells.Item(19, 4).Value = 735
$ws.Cells.Item(20, 3).Value = 1046
$ws.Cells.Item(20, 4).Value = 404
$ws.Cells.Item(21, 3).Value = 579
$ws.Cells.Item(21, 4).Value = 588
$ws.Cells.Item(22, 3).Value = 896
$ws.Cells.Item(22, 4).Value = 748
$ws.Cells.Item(23, 3).Value = 1038
$ws.Cells.Item(23, 4).Value = 472
$ws.Cells.Item(24, 3).Value = 274
$ws.Cells.Item(24, 4).Value = 622
$ws.Cells.Item(25, 3).Value = 59
$ws.Cells.Item(25, 4).Value = 268
$ws.Cells.Item(26, 3).Value = 85
$ws.Cells.Item(26, 4).Value = 453
$ws.Cells.Item(27, 3).Value = 547
$ws.Cells.Item(27, 4).Value = 738
$ws.Cells.Item(28, 3).Value = 942
$ws.Cells.Item(28, 4).Value = 330
$ws.Cells.Item(29, 3).Value = 275
$ws.Cells.Item(29, 4).Value = 441
$ws.Cells.Item(30, 3).Value = 978
$ws.Cells.Item(30, 4).Value = 496
$ws.Cells.Item(31, 3).Value = 635
$ws.Cells.Item(31, 4).Value = 469
$ws.Cells.Item(32, 3).Value = 392
$ws.Cells.Item(32, 4).Value = 227
$ws.Cells.Item(33, 3).Value = 680
$ws.Cells.Item(33, 4).Value = 380
$ws.Cells.Item(34, 3).Value = 154
$ws.Cells.Item(34, 4).Value = 206
